$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended below the existing table (rows 26-31).
# Column layout matches the header row:
# A=Date, B=vechical name, C=Shift, D=s.hour, E=c.hour, F=Rent, G=R.Hour,
# H=D.rate, I=d.quantity, J=earning, K=M.regular, L=M.Rate, M=Total, N=Description

$rows = @(
    @{ Row=26; A="21-02-2018"; B="v2"; C="/Night";    D=1500.0; E=1516.0; F=2000.0; G=16.0; H=66.0; I=150.0; J=32000.0; K="Engine oil `t        250"; L=1200.0; M=20900.0; N="nothng" },
    @{ Row=27; A="14-02-2018"; B="v3"; C="Day/";       D=1500.0; E=1516.0; F=1700.0; G=16.0; H=66.0; I=200.0; J=27200.0; K="Engine oil `t        250"; L=1200.0; M=12800.0; N="nothig" },
    @{ Row=28; A="26-02-2018"; B="v1"; C="Day/Night"; D=1500.0; E=1516.0; F=1700.0; G=16.0; H=66.0; I=200.0; J=27200.0; K="Engine oil `t        250"; L=1200.0; M=12800.0; N="nothing" },
    @{ Row=29; A="27-02-2018"; B="v1"; C="Day/";       D=1516.0; E=1532.0; F=1700.0; G=16.0; H=66.0; I=200.0; J=27200.0; K="Hydraulic oil`t        1000"; L=500.0;  M=13500.0; N="AAAA" },
    @{ Row=30; A="25-02-2018"; B="v1"; C="/Night";    D=1560.0; E=1572.0; F=1700.0; G=12.0; H=66.0; I=100.0; J=20400.0; K="Nothing"; L=200.0;  M=13600.0; N=",.JN,.J" },
    @{ Row=31; A="14-03-2018"; B="v1"; C="Day/Night"; D=1600.0; E=1616.0; F=1700.0; G=16.0; H=66.0; I=200.0; J=27200.0; K="Engine oil `t        250"; L=1200.0; M=12800.0; N="no" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
}
